$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered list of row -> "Exclusions after 1st run" value to add in column H
# (2017 end-of-season exclusion ranges carried over from the MacFerrin et
# al., 2019 outputs; 2018 still to be done per the commit message)
$updates = [ordered]@{
    80  = "3910-"
    82  = "1270-"
    86  = "760-900 1155-1315 3130-3170 3250-3300"
    87  = "1690-1900 2750-2815 3000-3070 3580-"
    89  = "855-945 2570-2650"
    90  = "2700-2935"
    91  = "-1645 1770-1865 2865-2950 2970-"
    92  = "535-1240 3090-3230"
    93  = "5090-6365"
    97  = "2700-2850"
    100 = "4145-"
    107 = "1740-2040 2870-3930"
    109 = "5255-5410 15660-15790 16495-"
    115 = "21700-21850"
    117 = "-315 950-1055"
    119 = "-200 660-800 1560-1575"
    123 = "530-630"
    127 = "-475 1130-1170 5875-5950 6190-6340 6715-7370 7460-7780 7900-"
    133 = "-175 290-340 630-970 1090-1220"
    134 = "1720-"
    136 = "-240 1085-1315 1390-1430 1450-1530 1740-1810 2030-2060"
    137 = "3770-3900 4265-4335"
    138 = "-70 2050-2310"
    140 = "6690-6730"
    141 = "6480-6505 7200-7250"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 8).Value = $updates[$row]
}

# Adjust column E width (auto-resized as a side effect of data entry)
$ws.Columns.Item(5).ColumnWidth = 20

# Update the view state to match the final saved selection
$ws.Application.ActiveWindow.ScrollRow = 125
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("G137").Select()
